$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Kompetenzraster")
$ws2 = $wb.Worksheets.Item("Berufe")
$ws = $ws2

$ws.Cells.Item(1,5).Value = "AgrarpraktikerIn EBA"
$ws.Cells.Item(2,5).Value = "agrarpra"
$ws.Cells.Item(1,6).Value = "DetailhandelsfachFrauMann EFZ"
$ws.Cells.Item(2,6).Value = "detailhandelsfach"
$ws.Cells.Item(1,8).Value = "alle As"
$ws.Cells.Item(1,7).Value = "InformatikerIn EFZ"
$ws.Cells.Item(2,7).Value = "informatik"
$ws.Cells.Item(3,5).Value = "x"
$ws.Cells.Item(3,6).Value = "x"
$ws.Cells.Item(3,7).Value = "x"
$ws.Cells.Item(3,8).Value = "x"
$ws.Cells.Item(4,5).Value = "x"
$ws.Cells.Item(4,6).Value = "x"
$ws.Cells.Item(4,7).Value = "x"
$ws.Cells.Item(4,8).Value = "x"
$ws.Cells.Item(5,5).Value = "x"
$ws.Cells.Item(5,6).Value = "x"
$ws.Cells.Item(5,7).Value = "x"
$ws.Cells.Item(6,5).Value = "x"
$ws.Cells.Item(6,6).Value = "x"
$ws.Cells.Item(6,7).Value = "x"
$ws.Cells.Item(7,7).Value = "x"
$ws.Cells.Item(8,7).Value = "x"
$ws.Cells.Item(9,7).Value = "x"
$ws.Cells.Item(10,5).Value = "x"
$ws.Cells.Item(10,6).Value = "x"
$ws.Cells.Item(10,7).Value = "x"
$ws.Cells.Item(10,8).Value = "x"
$ws.Cells.Item(11,5).Value = "x"
$ws.Cells.Item(11,7).Value = "x"
$ws.Cells.Item(12,7).Value = "x"
$ws.Cells.Item(13,7).Value = "x"
$ws.Cells.Item(14,7).Value = "x"
$ws.Cells.Item(15,5).Value = "x"
$ws.Cells.Item(15,6).Value = "x"
$ws.Cells.Item(15,7).Value = "x"
$ws.Cells.Item(15,8).Value = "x"
$ws.Cells.Item(16,5).Value = "x"
$ws.Cells.Item(16,6).Value = "x"
$ws.Cells.Item(16,7).Value = "x"
$ws.Cells.Item(16,8).Value = "x"
$ws.Cells.Item(17,5).Value = "x"
$ws.Cells.Item(17,7).Value = "x"
$ws.Cells.Item(18,7).Value = "x"
$ws.Cells.Item(19,7).Value = "x"
$ws.Cells.Item(20,5).Value = "x"
$ws.Cells.Item(20,6).Value = "x"
$ws.Cells.Item(20,7).Value = "x"
$ws.Cells.Item(20,8).Value = "x"
$ws.Cells.Item(21,5).Value = "x"
$ws.Cells.Item(21,6).Value = "x"
$ws.Cells.Item(21,7).Value = "x"
$ws.Cells.Item(22,5).Value = "x"
$ws.Cells.Item(22,7).Value = "x"
$ws.Cells.Item(23,5).Value = "x"
$ws.Cells.Item(23,7).Value = "x"
$ws.Cells.Item(24,7).Value = "x"
$ws.Cells.Item(25,5).Value = "x"
$ws.Cells.Item(25,6).Value = "x"
$ws.Cells.Item(25,7).Value = "x"
$ws.Cells.Item(25,8).Value = "x"
$ws.Cells.Item(26,5).Value = "x"
$ws.Cells.Item(26,6).Value = "x"
$ws.Cells.Item(26,7).Value = "x"
$ws.Cells.Item(27,6).Value = "x"
$ws.Cells.Item(27,7).Value = "x"
$ws.Cells.Item(28,7).Value = "x"
$ws.Cells.Item(29,5).Value = "x"
$ws.Cells.Item(29,6).Value = "x"
$ws.Cells.Item(29,7).Value = "x"
$ws.Cells.Item(29,8).Value = "x"
$ws.Cells.Item(30,7).Value = "x"
$ws.Cells.Item(31,7).Value = "x"
$ws.Cells.Item(32,7).Value = "x"
$ws.Cells.Item(33,7).Value = "x"
$ws.Cells.Item(34,5).Value = "x"
$ws.Cells.Item(34,6).Value = "x"
$ws.Cells.Item(34,7).Value = "x"
$ws.Cells.Item(34,8).Value = "x"
$ws.Cells.Item(35,6).Value = "x"
$ws.Cells.Item(35,7).Value = "x"
$ws.Cells.Item(36,7).Value = "x"
$ws.Cells.Item(37,7).Value = "x"
$ws.Cells.Item(38,5).Value = "x"
$ws.Cells.Item(38,6).Value = "x"
$ws.Cells.Item(38,7).Value = "x"
$ws.Cells.Item(38,8).Value = "x"
$ws.Cells.Item(39,6).Value = "x"
$ws.Cells.Item(39,7).Value = "x"
$ws.Cells.Item(40,7).Value = "x"
$ws.Cells.Item(41,7).Value = "x"
$ws.Cells.Item(42,7).Value = "x"
$ws.Cells.Item(43,5).Value = "x"
$ws.Cells.Item(43,6).Value = "x"
$ws.Cells.Item(43,7).Value = "x"
$ws.Cells.Item(43,8).Value = "x"
$ws.Cells.Item(44,7).Value = "x"
$ws.Cells.Item(45,7).Value = "x"
$ws.Cells.Item(46,7).Value = "x"
$ws.Cells.Item(47,7).Value = "x"
$ws.Cells.Item(48,7).Value = "x"
$ws.Cells.Item(49,5).Value = "x"
$ws.Cells.Item(49,6).Value = "x"
$ws.Cells.Item(49,7).Value = "x"
$ws.Cells.Item(49,8).Value = "x"
$ws.Cells.Item(52,7).Value = "x"
$ws.Cells.Item(53,7).Value = "x"
$ws.Cells.Item(54,5).Value = "x"
$ws.Cells.Item(54,6).Value = "x"
$ws.Cells.Item(54,7).Value = "x"
$ws.Cells.Item(54,8).Value = "x"
$ws.Cells.Item(55,7).Value = "x"
$ws.Cells.Item(56,7).Value = "x"
$ws.Cells.Item(57,7).Value = "x"
$ws.Cells.Item(58,7).Value = "x"
$ws.Cells.Item(59,7).Value = "x"
$ws.Cells.Item(61,5).Value = "x"
$ws.Cells.Item(61,6).Value = "x"
$ws.Cells.Item(61,7).Value = "x"
$ws.Cells.Item(61,8).Value = "x"
$ws.Cells.Item(62,5).Value = "x"
$ws.Cells.Item(64,7).Value = "x"
$ws.Cells.Item(65,7).Value = "x"
$ws.Cells.Item(66,7).Value = "x"
$ws.Cells.Item(67,7).Value = "x"
$ws.Cells.Item(68,5).Value = "x"
$ws.Cells.Item(68,6).Value = "x"
$ws.Cells.Item(68,7).Value = "x"
$ws.Cells.Item(68,8).Value = "x"
$ws.Cells.Item(74,5).Value = "x"
$ws.Cells.Item(74,6).Value = "x"
$ws.Cells.Item(74,7).Value = "x"
$ws.Cells.Item(74,8).Value = "x"
$ws.Cells.Item(75,7).Value = "x"
$ws.Cells.Item(77,7).Value = "x"
$ws.Cells.Item(80,5).Value = "x"
$ws.Cells.Item(80,6).Value = "x"
$ws.Cells.Item(80,7).Value = "x"
$ws.Cells.Item(80,8).Value = "x"
$ws.Cells.Item(85,5).Value = "x"
$ws.Cells.Item(85,6).Value = "x"
$ws.Cells.Item(85,7).Value = "x"
$ws.Cells.Item(85,8).Value = "x"
$ws.Cells.Item(89,5).Value = "x"
$ws.Cells.Item(89,6).Value = "x"
$ws.Cells.Item(89,7).Value = "x"
$ws.Cells.Item(89,8).Value = "x"
$ws.Cells.Item(94,5).Value = "x"
$ws.Cells.Item(94,6).Value = "x"
$ws.Cells.Item(94,7).Value = "x"
$ws.Cells.Item(94,8).Value = "x"
$ws.Cells.Item(95,5).Value = "x"
$ws.Cells.Item(95,7).Value = "x"
$ws.Cells.Item(96,7).Value = "x"
$ws.Cells.Item(97,7).Value = "x"
$ws.Cells.Item(98,5).Value = "x"
$ws.Cells.Item(98,6).Value = "x"
$ws.Cells.Item(98,7).Value = "x"
$ws.Cells.Item(98,8).Value = "x"
$ws.Cells.Item(99,5).Value = "x"
$ws.Cells.Item(99,6).Value = "x"
$ws.Cells.Item(99,7).Value = "x"
$ws.Cells.Item(100,7).Value = "x"
$ws.Cells.Item(101,7).Value = "x"
$ws.Cells.Item(102,7).Value = "x"
$ws.Cells.Item(103,7).Value = "x"
$ws.Cells.Item(104,5).Value = "x"
$ws.Cells.Item(104,6).Value = "x"
$ws.Cells.Item(104,7).Value = "x"
$ws.Cells.Item(104,8).Value = "x"
$ws.Cells.Item(105,7).Value = "x"
$ws.Cells.Item(106,7).Value = "x"
$ws.Cells.Item(107,7).Value = "x"
$ws.Cells.Item(108,5).Value = "x"
$ws.Cells.Item(108,6).Value = "x"
$ws.Cells.Item(108,7).Value = "x"
$ws.Cells.Item(108,8).Value = "x"
$ws.Cells.Item(109,5).Value = "x"
$ws.Cells.Item(109,6).Value = "x"
$ws.Cells.Item(109,7).Value = "x"
$ws.Cells.Item(110,5).Value = "x"
$ws.Cells.Item(110,6).Value = "x"
$ws.Cells.Item(110,7).Value = "x"
$ws.Cells.Item(111,7).Value = "x"
$ws.Cells.Item(112,7).Value = "x"
$ws.Cells.Item(113,5).Value = "x"
$ws.Cells.Item(113,6).Value = "x"
$ws.Cells.Item(113,7).Value = "x"
$ws.Cells.Item(113,8).Value = "x"
$ws.Cells.Item(114,5).Value = "x"
$ws.Cells.Item(114,6).Value = "x"
$ws.Cells.Item(114,7).Value = "x"
$ws.Cells.Item(115,7).Value = "x"
$ws.Cells.Item(116,7).Value = "x"
$ws.Cells.Item(117,7).Value = "x"
$ws.Cells.Item(118,7).Value = "x"
$ws.Cells.Item(119,7).Value = "x"
$ws.Cells.Item(120,7).Value = "x"
$ws.Cells.Item(121,7).Value = "x"
$ws.Cells.Item(122,7).Value = "x"
$ws.Cells.Item(123,7).Value = "x"
$ws.Cells.Item(124,7).Value = "x"

$ws.Columns.Item(1).ColumnWidth = 80.83072916666667
$ws.Columns.Item(6).ColumnWidth = 23.498697916666668

$ws1.Range("G7").Select()
$ws2.Activate()
$ws2.Range("G125").Select()
